$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 444.94287
$ws.Range("I19").Value = 411.05884
$ws.Range("J19").Value = 476.94446
$ws.Range("K19").Value = 411.05884
$ws.Range("L19").Value = 476.94446
$ws.Range("M19").Value = -236.05884
$ws.Range("N19").Value = -826.9444599999999
$ws.Range("H33").Value = 2066876
$ws.Range("I33").Value = 763.0857
$ws.Range("K33").Value = 763.0857
$ws.Range("M33").Value = -534.0857
$ws.Range("H40").Value = 11112978
$ws.Range("I40").Value = 1861.8387
$ws.Range("K40").Value = 1861.8387
$ws.Range("M40").Value = -1686.8387
$ws.Range("H64").Value = 4594.737
$ws.Range("I64").Value = 5740
$ws.Range("K64").Value = 5740
$ws.Range("M64").Value = -5492
$ws.Range("H67").Value = 4594.737
$ws.Range("I67").Value = 5740
$ws.Range("K67").Value = 5740
$ws.Range("M67").Value = -4882
$ws.Range("H70").Value = 4750.25
$ws.Range("I70").Value = 6001
$ws.Range("J70").Value = 4333.3335
$ws.Range("K70").Value = 18003
$ws.Range("L70").Value = 13000.0005
$ws.Range("M70").Value = -17733
$ws.Range("N70").Value = -13540.0005
$ws.Range("H73").Value = 4750.25
$ws.Range("I73").Value = 6001
$ws.Range("J73").Value = 4333.3335
$ws.Range("K73").Value = 18003
$ws.Range("L73").Value = 13000.0005
$ws.Range("M73").Value = -17067
$ws.Range("N73").Value = -14872.0005
$ws.Range("H115").Value = 385
$ws.Range("I115").Value = 385
$ws.Range("K115").Value = 1155
$ws.Range("M115").Value = 412
$ws.Range("H116").Value = 8854.556
$ws.Range("I116").Value = 9240.117
$ws.Range("K116").Value = 9240.117
$ws.Range("M116").Value = -5798.117
$ws.Range("H132").Value = 4168437.8
$ws.Range("I132").Value = 1405.4203
$ws.Range("K132").Value = 4216.2609
$ws.Range("M132").Value = -1686.2609
$ws.Range("H137").Value = 1331.9429
$ws.Range("I137").Value = 1391.5238
$ws.Range("K137").Value = 4174.5714
$ws.Range("M137").Value = -1624.5714

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 142859020
$ws.Range("I63").Value = 166668600
$ws.Range("K63").Value = 166668600
$ws.Range("M63").Value = -166667914
$ws.Range("H66").Value = 142859020
$ws.Range("I66").Value = 166668600
$ws.Range("K66").Value = 833343000
$ws.Range("M66").Value = -833339568
$ws.Range("H110").Value = 7425.868
$ws.Range("J110").Value = 2143.3635
$ws.Range("L110").Value = 2143.3635
$ws.Range("N110").Value = -6233.363499999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H99").Value = 62502170
$ws.Range("I99").Value = 83335144
$ws.Range("J99").Value = 3250
$ws.Range("K99").Value = 83335144
$ws.Range("L99").Value = 3250
$ws.Range("M99").Value = -83333646
$ws.Range("N99").Value = -6246
$ws.Range("H134").Value = 204425.2
$ws.Range("I134").Value = 5528
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 16584
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -14049
$ws.Range("N134").Value = -3005112

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9441013
$ws.Range("I31").Value = 1719.6842
$ws.Range("J31").Value = 14715913
$ws.Range("K31").Value = 1719.6842
$ws.Range("L31").Value = 14715913
$ws.Range("M31").Value = -1424.6842
$ws.Range("N31").Value = -14716503
$ws.Range("H34").Value = 9441013
$ws.Range("I34").Value = 1719.6842
$ws.Range("J34").Value = 14715913
$ws.Range("K34").Value = 1719.6842
$ws.Range("L34").Value = 14715913
$ws.Range("M34").Value = -1517.6842
$ws.Range("N34").Value = -14716317
$ws.Range("H86").Value = 1398.75
$ws.Range("I86").Value = 1342.1428
$ws.Range("J86").Value = 1478
$ws.Range("K86").Value = 1342.1428
$ws.Range("L86").Value = 1478
$ws.Range("M86").Value = -219.1428000000001
$ws.Range("N86").Value = -3724
$ws.Range("H89").Value = 1398.75
$ws.Range("I89").Value = 1342.1428
$ws.Range("J89").Value = 1478
$ws.Range("K89").Value = 6710.714
$ws.Range("L89").Value = 7390
$ws.Range("M89").Value = -1094.714
$ws.Range("N89").Value = -18622
$ws.Range("H99").Value = 10418958
$ws.Range("I99").Value = 2080
$ws.Range("J99").Value = 17859586
$ws.Range("K99").Value = 2080
$ws.Range("L99").Value = 17859586
$ws.Range("M99").Value = -582
$ws.Range("N99").Value = -17862582
$ws.Range("H126").Value = 10418958
$ws.Range("I126").Value = 2080
$ws.Range("J126").Value = 17859586
$ws.Range("K126").Value = 6240
$ws.Range("L126").Value = 53578758
$ws.Range("M126").Value = -3770
$ws.Range("N126").Value = -53583698

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4615875.5
$ws.Range("I113").Value = 50000000
$ws.Range("J113").Value = 833865.3
$ws.Range("K113").Value = 150000000
$ws.Range("L113").Value = 2501595.9
$ws.Range("M113").Value = -149997830
$ws.Range("N113").Value = -2505935.9
$ws.Range("H121").Value = 761.6
$ws.Range("I121").Value = 269.375
$ws.Range("J121").Value = 940.5909
$ws.Range("K121").Value = 808.125
$ws.Range("L121").Value = 2821.7727
$ws.Range("M121").Value = 501.875
$ws.Range("N121").Value = -5441.7727
$ws.Range("H131").Value = 2128652.8
$ws.Range("I131").Value = 6250660.5
$ws.Range("J131").Value = 1164.871
$ws.Range("K131").Value = 18751981.5
$ws.Range("L131").Value = 3494.613
$ws.Range("M131").Value = -18746941.5
$ws.Range("N131").Value = -13574.613

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 6791.125
$ws.Range("J57").Value = 6791.125
$ws.Range("L57").Value = 6791.125
$ws.Range("N57").Value = -8431.125
$ws.Range("H107").Value = 890.05884
$ws.Range("I107").Value = 198.63637
$ws.Range("J107").Value = 2157.6667
$ws.Range("K107").Value = 198.63637
$ws.Range("L107").Value = 2157.6667
$ws.Range("M107").Value = 1721.36363
$ws.Range("N107").Value = -5997.6667

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3627
$ws.Range("I7").Value = 2143.4285
$ws.Range("K7").Value = 2143.4285
$ws.Range("M7").Value = -2031.4285
$ws.Range("H126").Value = 3627
$ws.Range("I126").Value = 2143.4285
$ws.Range("K126").Value = 6430.2855
$ws.Range("M126").Value = -3960.2855

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 469444
$ws.Range("I62").Value = 469444
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 469444
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -468820
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 469444
$ws.Range("I65").Value = 469444
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 2347220
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -2344100
$ws.Range("N65").ClearContents()
$ws.Range("H136").Value = 5956253
$ws.Range("I136").Value = 5052.6924
$ws.Range("K136").Value = 15158.0772
$ws.Range("M136").Value = -12608.0772
